$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Update columns D (Update Test Passed) and E (Delete Test Passed) to TRUE
# for rows 2-15 and 20-24 (rows 16-19 already TRUE and unchanged).
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,20,21,22,23,24)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = $true
}

# Update the active selection on the Test Results sheet
$ws.Activate()
$ws.Range("G21").Select()
